$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Replace("0x7f776b197a60", "0x7f98e35edc70") | Out-Null
$ws.Range("C2").Replace("0x7f776b0e8880", "0x7f8fac0af3d0") | Out-Null

$ws.Range("B3").Replace("0x7f78564dcb20", "0x7f98e0da3250") | Out-Null
$ws.Range("C3").Replace("0x7f776b0f68b0", "0x7f9864250f40") | Out-Null

$ws.Range("B4").Replace("0x7f77640a35e0", "0x7f8fac269d00") | Out-Null
$ws.Range("C4").Replace("0x7f776b6d0640", "0x7f90f8024af0") | Out-Null

$ws.Range("B5").Replace("0x7f77640a37f0", "0x7f991774a280") | Out-Null
$ws.Range("C5").Replace("0x7f7763cdf4c0", "0x7f8fac0c4b20") | Out-Null

$ws.Range("B6").Replace("0x7f7764080190", "0x7f90f80d3e20") | Out-Null
$ws.Range("C6").Replace("0x7f7763cf6d90", "0x7f8f6c6adca0") | Out-Null
